$d = $word.ActiveDocument

# 1) Dateline: "Valparaíso, 24" -> "Santiago, 28" (scope to the dateline
#    paragraph so we don't touch "Universidad de Valparaíso" later on).
$pDate = $d.Paragraphs(2).Range
$pDate.Find.Execute("Valparaíso", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Santiago", 2)

$pDate2 = $d.Paragraphs(2).Range
$pDate2.Find.Execute(", 24", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", 28", 2)

# 2) Signatory name
$d.Content.Find.Execute("Eduardo Muñoz Inchausti", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Jeannette Rodríguez Chandia", 2)

# 3) Signatory title
$d.Content.Find.Execute("Director", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Jefa de Carrera Campus Santiago", 2)

# 4) Initials line "EMI/krr" -> "JRC/lll"
$d.Content.Find.Execute("EMI", $true, $false, $false, $false, $false, `
    $true, 1, $false, "JRC", 2)

$d.Content.Find.Execute("krr", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lll", 2)

# 5) Footer address/phone line
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute( `
    "Las Heras Nº 06 Valparaíso | Fono: (32) 250 7961- 2507815 | E-mail: practivasv@uv.cl, www.uv.cl", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Campus Santiago - Gran Avenida 4160, San Miguel | Fono +56 (2)2329  2149", 2)
